$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows ---
# O1485: 0 -> 2
$ws.Range("O1485").Value = 2

# R1487 and R1488: inlineStr empty -> numeric 0
$ws.Range("R1487").Value = 0
$ws.Range("R1488").Value = 0

# --- Append new weekly rows 1489-1509 (A:Q); column R intentionally left blank to match source (empty cell) ---
$data = New-Object 'object[,]' 21,17
$data[0,0] = 45474
$data[0,1] = 1559.5
$data[0,2] = 1665.849975585938
$data[0,3] = 1559.5
$data[0,4] = 1647.449951171875
$data[0,5] = 1628.903076171875
$data[0,6] = 39638242
$data[0,7] = 2024
$data[0,8] = 7
$data[0,9] = 1
$data[0,10] = 0
$data[0,11] = 0
$data[0,12] = 0
$data[0,13] = 27
$data[0,14] = 0
$data[0,15] = 0
$data[0,16] = 0
$data[1,0] = 45481
$data[1,1] = 1643.099975585938
$data[1,2] = 1719.75
$data[1,3] = 1637.550048828125
$data[1,4] = 1711.75
$data[1,5] = 1692.479248046875
$data[1,6] = 46860480
$data[1,7] = 2024
$data[1,8] = 7
$data[1,9] = 8
$data[1,10] = 0
$data[1,11] = 0
$data[1,12] = 0
$data[1,13] = 28
$data[1,14] = 0
$data[1,15] = 0
$data[1,16] = 0
$data[2,0] = 45488
$data[2,1] = 1726.199951171875
$data[2,2] = 1844
$data[2,3] = 1700
$data[2,4] = 1792.949951171875
$data[2,5] = 1772.765014648438
$data[2,6] = 54896523
$data[2,7] = 2024
$data[2,8] = 7
$data[2,9] = 15
$data[2,10] = 0
$data[2,11] = 0
$data[2,12] = 0
$data[2,13] = 29
$data[2,14] = 0
$data[2,15] = 1
$data[2,16] = 1
$data[3,0] = 45495
$data[3,1] = 1800.150024414062
$data[3,2] = 1883
$data[3,3] = 1783.25
$data[3,4] = 1878.900024414062
$data[3,5] = 1857.747436523438
$data[3,6] = 38762673
$data[3,7] = 2024
$data[3,8] = 7
$data[3,9] = 22
$data[3,10] = 0
$data[3,11] = 0
$data[3,12] = 0
$data[3,13] = 30
$data[3,14] = 0
$data[3,15] = 0
$data[3,16] = 0
$data[4,0] = 45502
$data[4,1] = 1893.599975585938
$data[4,2] = 1903
$data[4,3] = 1816.150024414062
$data[4,4] = 1821.199951171875
$data[4,5] = 1800.697021484375
$data[4,6] = 27237357
$data[4,7] = 2024
$data[4,8] = 7
$data[4,9] = 29
$data[4,10] = 0
$data[4,11] = 0
$data[4,12] = 0
$data[4,13] = 31
$data[4,14] = 0
$data[4,15] = 0
$data[4,16] = 0
$data[5,0] = 45509
$data[5,1] = 1784.949951171875
$data[5,2] = 1797.900024414062
$data[5,3] = 1718.550048828125
$data[5,4] = 1770.75
$data[5,5] = 1750.81494140625
$data[5,6] = 34007629
$data[5,7] = 2024
$data[5,8] = 8
$data[5,9] = 5
$data[5,10] = 0
$data[5,11] = 0
$data[5,12] = 0
$data[5,13] = 32
$data[5,14] = 0
$data[5,15] = 0
$data[5,16] = 0
$data[6,0] = 45516
$data[6,1] = 1773.050048828125
$data[6,2] = 1861.849975585938
$data[6,3] = 1768.099975585938
$data[6,4] = 1858.949951171875
$data[6,5] = 1838.02197265625
$data[6,6] = 23278622
$data[6,7] = 2024
$data[6,8] = 8
$data[6,9] = 12
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 0
$data[6,13] = 33
$data[6,14] = 0
$data[6,15] = 0
$data[6,16] = 0
$data[7,0] = 45523
$data[7,1] = 1860
$data[7,2] = 1893.650024414062
$data[7,3] = 1849.199951171875
$data[7,4] = 1862.099975585938
$data[7,5] = 1841.136596679688
$data[7,6] = 17739446
$data[7,7] = 2024
$data[7,8] = 8
$data[7,9] = 19
$data[7,10] = 0
$data[7,11] = 0
$data[7,12] = 0
$data[7,13] = 34
$data[7,14] = 0
$data[7,15] = 0
$data[7,16] = 0
$data[8,0] = 45530
$data[8,1] = 1870
$data[8,2] = 1951
$data[8,3] = 1864
$data[8,4] = 1943.699951171875
$data[8,5] = 1921.81787109375
$data[8,6] = 34768842
$data[8,7] = 2024
$data[8,8] = 8
$data[8,9] = 26
$data[8,10] = 0
$data[8,11] = 0
$data[8,12] = 0
$data[8,13] = 35
$data[8,14] = 0
$data[8,15] = 0
$data[8,16] = 0
$data[9,0] = 45537
$data[9,1] = 1943.349975585938
$data[9,2] = 1975.75
$data[9,3] = 1896.849975585938
$data[9,4] = 1901.849975585938
$data[9,5] = 1880.439086914062
$data[9,6] = 24976287
$data[9,7] = 2024
$data[9,8] = 9
$data[9,9] = 2
$data[9,10] = 0
$data[9,11] = 0
$data[9,12] = 0
$data[9,13] = 36
$data[9,14] = 0
$data[9,15] = 0
$data[9,16] = 0
$data[10,0] = 45544
$data[10,1] = 1890
$data[10,2] = 1958.599975585938
$data[10,3] = 1889
$data[10,4] = 1944.099975585938
$data[10,5] = 1922.21337890625
$data[10,6] = 29393377
$data[10,7] = 2024
$data[10,8] = 9
$data[10,9] = 9
$data[10,10] = 0
$data[10,11] = 0
$data[10,12] = 0
$data[10,13] = 37
$data[10,14] = 0
$data[10,15] = 0
$data[10,16] = 0
$data[11,0] = 45551
$data[11,1] = 1945.75
$data[11,2] = 1958.449951171875
$data[11,3] = 1867.400024414062
$data[11,4] = 1905.75
$data[11,5] = 1884.295166015625
$data[11,6] = 33968320
$data[11,7] = 2024
$data[11,8] = 9
$data[11,9] = 16
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 0
$data[11,13] = 38
$data[11,14] = 0
$data[11,15] = 0
$data[11,16] = 0
$data[12,0] = 45558
$data[12,1] = 1909
$data[12,2] = 1974.599975585938
$data[12,3] = 1871.400024414062
$data[12,4] = 1906.75
$data[12,5] = 1885.283935546875
$data[12,6] = 39044363
$data[12,7] = 2024
$data[12,8] = 9
$data[12,9] = 23
$data[12,10] = 0
$data[12,11] = 0
$data[12,12] = 0
$data[12,13] = 39
$data[12,14] = 0
$data[12,15] = 0
$data[12,16] = 0
$data[13,0] = 45565
$data[13,1] = 1880
$data[13,2] = 1954.099975585938
$data[13,3] = 1870.5
$data[13,4] = 1918.150024414062
$data[13,5] = 1896.555541992188
$data[13,6] = 28521397
$data[13,7] = 2024
$data[13,8] = 9
$data[13,9] = 30
$data[13,10] = 0
$data[13,11] = 0
$data[13,12] = 0
$data[13,13] = 40
$data[13,14] = 0
$data[13,15] = 0
$data[13,16] = 0
$data[14,0] = 45572
$data[14,1] = 1923.199951171875
$data[14,2] = 1977
$data[14,3] = 1906.349975585938
$data[14,4] = 1935.099975585938
$data[14,5] = 1913.314697265625
$data[14,6] = 26098679
$data[14,7] = 2024
$data[14,8] = 10
$data[14,9] = 7
$data[14,10] = 0
$data[14,11] = 0
$data[14,12] = 0
$data[14,13] = 41
$data[14,14] = 0
$data[14,15] = 0
$data[14,16] = 0
$data[15,0] = 45579
$data[15,1] = 1932.150024414062
$data[15,2] = 1991.449951171875
$data[15,3] = 1869.25
$data[15,4] = 1879.599975585938
$data[15,5] = 1858.439575195312
$data[15,6] = 30905613
$data[15,7] = 2024
$data[15,8] = 10
$data[15,9] = 14
$data[15,10] = 0
$data[15,11] = 0
$data[15,12] = 0
$data[15,13] = 42
$data[15,14] = 1
$data[15,15] = 0
$data[15,16] = 0
$data[16,0] = 45586
$data[16,1] = 1891.099975585938
$data[16,2] = 1898.800048828125
$data[16,3] = 1838
$data[16,4] = 1862.050048828125
$data[16,5] = 1841.087158203125
$data[16,6] = 20687346
$data[16,7] = 2024
$data[16,8] = 10
$data[16,9] = 21
$data[16,10] = 0
$data[16,11] = 0
$data[16,12] = 0
$data[16,13] = 43
$data[16,14] = 0
$data[16,15] = 0
$data[16,16] = 0
$data[17,0] = 45593
$data[17,1] = 1859.699951171875
$data[17,2] = 1881.900024414062
$data[17,3] = 1746.5
$data[17,4] = 1760.849975585938
$data[17,5] = 1741.0263671875
$data[17,6] = 25327109
$data[17,7] = 2024
$data[17,8] = 10
$data[17,9] = 28
$data[17,10] = 0
$data[17,11] = 0
$data[17,12] = 0
$data[17,13] = 44
$data[17,14] = 0
$data[17,15] = 0
$data[17,16] = 0
$data[18,0] = 45600
$data[18,1] = 1724.5
$data[18,2] = 1840.599975585938
$data[18,3] = 1718
$data[18,4] = 1829.949951171875
$data[18,5] = 1829.949951171875
$data[18,6] = 29681021
$data[18,7] = 2024
$data[18,8] = 11
$data[18,9] = 4
$data[18,10] = 0
$data[18,11] = 0
$data[18,12] = 0
$data[18,13] = 45
$data[18,14] = 0
$data[18,15] = 0
$data[18,16] = 0
$data[19,0] = 45607
$data[19,1] = 1829
$data[19,2] = 1881
$data[19,3] = 1822.550048828125
$data[19,4] = 1864.550048828125
$data[19,5] = 1864.550048828125
$data[19,6] = 17385610
$data[19,7] = 2024
$data[19,8] = 11
$data[19,9] = 11
$data[19,10] = 0
$data[19,11] = 0
$data[19,12] = 0
$data[19,13] = 46
$data[19,14] = 0
$data[19,15] = 0
$data[19,16] = 0
$data[20,0] = 45614
$data[20,1] = 1849.199951171875
$data[20,2] = 1914.050048828125
$data[20,3] = 1795
$data[20,4] = 1902.25
$data[20,5] = 1902.25
$data[20,6] = 21245150
$data[20,7] = 2024
$data[20,8] = 11
$data[20,9] = 18
$data[20,10] = 0
$data[20,11] = 0
$data[20,12] = 0
$data[20,13] = 47
$data[20,14] = 0
$data[20,15] = 0
$data[20,16] = 0
$ws.Range("A1489:Q1509").Value = $data
$ws.Range("A1489:A1509").NumberFormat = "YYYY-MM-DD HH:MM:SS"
